$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on affected Price/Volume columns so numeric-looking
# strings (and percentages) are stored as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.63%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.47%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.070'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.45%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05691'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.43%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.485'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.16%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8198'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.76%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8422'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.41%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1329'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.88%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06891'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.89%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02860'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.51%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09393'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.03%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001522'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.39%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04100'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-11.68%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005986'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.40%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006143'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.04%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '3,761.17%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.510'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.48%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.006'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.11%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.229'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '8.43%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.88%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1298'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.72%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.554'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-5.06%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.67%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.001216'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.82%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.003966'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-13.60%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.99%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.82%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005497'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-11.37%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1057'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.10%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002469'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.29%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009405'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '5.56%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005208'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.71%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.07%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1015'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-15.48%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002595'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3.21%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.07%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.07%'
